$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy styles from column M (existing) to new column N for rows 2-10
$ws.Range("N2").Style = $ws.Range("M2").Style
$ws.Range("N3").Style = $ws.Range("M3").Style
$ws.Range("N4").Style = $ws.Range("M4").Style
$ws.Range("N5").Style = $ws.Range("M5").Style
$ws.Range("N6").Style = $ws.Range("M6").Style
$ws.Range("N7").Style = $ws.Range("M7").Style
$ws.Range("N8").Style = $ws.Range("M8").Style
$ws.Range("N9").Style = $ws.Range("M9").Style
$ws.Range("N10").Style = $ws.Range("M10").Style

# New column N values
$ws.Range("N3").Value = 2022
$ws.Range("N4").Value = 1434
$ws.Range("N5").Value = 12822
$ws.Range("N6").Value = 3099
$ws.Range("N7").Value = 9722
$ws.Range("N8").Value = 14424
$ws.Range("N9").Value = 5279
$ws.Range("N10").Value = 9145

# Update selection to reflect the new active cell
$ws.Range("N2").Select()
